# This script re-shuffles the species-occurrence data among several rows of
# the "Artfynd" sheet. The rows below keep their location-independent columns
# (P, S, T, U, V, W, AW, AX, ...) fixed, but the observation-specific columns
# (Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn, Vetenskapligt
# namn, Auktor, Ost, Nord, Startdatum, Slutdatum) are cyclically rotated
# between the rows that make up each group below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together with an observation record.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Y", "AA")

# Columns that hold plain-text, date-shaped strings (e.g. "2023-08-26").
# These must be written back as TEXT, not auto-converted to a date serial
# number by Excel, so they need special handling when written.
$dateTextCols = @("Y", "AA")

# Groups of rows whose observation data rotate among themselves.
# new_data(row[i]) = old_data(row[i+1]), wrapping around the group.
$cycles = @(
    , @(9, 10)
    , @(18, 22, 23)
    , @(36, 38)
    , @(37, 39)
    , @(49, 55, 52, 53, 50, 56, 54, 51)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Count

    # Snapshot the current ("old") values for every tracked column of every
    # row in this cycle before any writes happen.
    $snapshot = @{}
    foreach ($row in $cycle) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$row").Value2
        }
        $snapshot[$row] = $rowVals
    }

    # Write the rotated values: row i receives the snapshot of row i+1.
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $n]
        $srcVals = $snapshot[$srcRow]
        foreach ($col in $cols) {
            $destCell = $ws.Range("$col$destRow")
            if ($dateTextCols -contains $col) {
                # Force text storage so "2023-08-26" stays a string instead
                # of becoming a date serial number.
                $destCell.NumberFormat = "@"
                $destCell.Value2 = $srcVals[$col]
                $destCell.Style = "Normal"
            } else {
                $destCell.Value2 = $srcVals[$col]
            }
        }
    }
}
